$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Auto fill: the quantity for the "Camera" row was filled/edited from 1 to 2.
# This cascades through the shared formula in E3 (=C3*D3) and the totals
# in D5 (=SUM(D2:D4)) and E5 automatically.
$ws.Range("D3").Value = 2

# Autosum: replace the manual addition in E5 with a proper SUM formula,
# mirroring the one already used in D5.
$ws.Range("E5").Formula = "=SUM(E2:E4)"

# Leave the selection where the user ended up after these edits.
$ws.Range("F11").Select()
